$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make sure the date cell for the new last row keeps the same date number format
# as all the other rows in column D (set before writing the value so Excel does
# not invent a brand new style entry for it).
$ws.Range("D113").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Shift existing data rows 33-112 down to 34-113 (values only; formats/styles already in place)
$src = $ws.Range("A33:R112")
$vals = $src.Value()
$dst = $ws.Range("A34:R113")
$dst.Value = $vals

# Write the new latest weekly record into row 33
$ws.Range("D33").Value = [DateTime]"2023-01-27"
$ws.Range("J33").Value = 50
$ws.Range("K33").Value = 40000
$ws.Range("L33").Value = 40000
$ws.Range("M33").Value = 40000
$ws.Range("O33").Value = "Región del Maule"
$ws.Range("P33").Value = 1600
